$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Target data for columns B:F, rows 2-50 (row 2 = ticker index 0) ---
# Each entry: row number => @(B, C, D, E, F)  ("" marks an intentionally blank cell)
$data = @{
    2 = @("NSE:AVALON", "NSE:3IINFOLTD", "NSE:MCX", "NSE:BANDHANBNK", "NSE:HDFCAMC")
    3 = @("NSE:COUNCODOS", "NSE:ABCAPITAL", "NSE:MUTHOOTFIN", "NSE:COALINDIA", "NSE:NAUKRI")
    4 = @("NSE:EMSLIMITED", "NSE:ADVENZYMES", "", "", "")
    5 = @("NSE:GULFOILLUB", "NSE:AGROPHOS", "", "", "")
    6 = @("NSE:HDFCAMC", "NSE:ANANDRATHI", "", "", "")
    7 = @("NSE:INFOMEDIA", "NSE:APOLLO", "", "", "")
    8 = @("NSE:JBMA", "NSE:ARIES", "", "", "")
    9 = @("NSE:KPIL", "NSE:ASAL", "", "", "")
    10 = @("NSE:NDGL", "NSE:ASKAUTOLTD", "", "", "")
    11 = @("NSE:NEULANDLAB", "NSE:BHAGYANGR", "", "", "")
    12 = @("NSE:NRBBEARING", "NSE:BHARATRAS", "", "", "")
    13 = @("NSE:PGIL", "NSE:BLS", "", "", "")
    14 = @("NSE:PIXTRANS", "NSE:BLUEJET", "", "", "")
    15 = @("NSE:PLASTIBLEN", "NSE:COROMANDEL", "", "", "")
    16 = @("NSE:PNBGILTS", "NSE:DALBHARAT", "", "", "")
    17 = @("NSE:PPAP", "NSE:DCAL", "", "", "")
    18 = @("NSE:QUESS", "NSE:DEEPAKNTR", "", "", "")
    19 = @("NSE:RITCO", "NSE:GNFC", "", "", "")
    20 = @("NSE:RRKABEL", "NSE:GSFC", "", "", "")
    21 = @("", "NSE:GSLSU", "", "", "")
    22 = @("", "NSE:GSPL", "", "", "")
    23 = @("", "NSE:GUJALKALI", "", "", "")
    24 = @("", "NSE:HINDPETRO", "", "", "")
    25 = @("", "NSE:INDOAMIN", "", "", "")
    26 = @("", "NSE:JBCHEPHARM", "", "", "")
    27 = @("", "NSE:KTKBANK", "", "", "")
    28 = @("", "NSE:LAOPALA", "", "", "")
    29 = @("", "NSE:LICI", "", "", "")
    30 = @("", "NSE:LINDEINDIA", "", "", "")
    31 = @("", "NSE:LTFOODS", "", "", "")
    32 = @("", "NSE:LUXIND", "", "", "")
    33 = @("", "NSE:MADRASFERT", "", "", "")
    34 = @("", "NSE:MANAKSIA", "", "", "")
    35 = @("", "NSE:MARATHON", "", "", "")
    36 = @("", "NSE:MOL", "", "", "")
    37 = @("", "NSE:NAGAFERT", "", "", "")
    38 = @("", "NSE:NFL", "", "", "")
    39 = @("", "NSE:ORIENTHOT", "", "", "")
    40 = @("", "NSE:PARADEEP", "", "", "")
    41 = @("", "NSE:PNC", "", "", "")
    42 = @("", "NSE:PRSMJOHNSN", "", "", "")
    43 = @("", "NSE:QUICKHEAL", "", "", "")
    44 = @("", "NSE:RAJESHEXPO", "", "", "")
    45 = @("", "NSE:RALLIS", "", "", "")
    46 = @("", "NSE:RAMASTEEL", "", "", "")
    47 = @("", "NSE:RCF", "", "", "")
    48 = @("", "NSE:RKEC", "", "", "")
    49 = @("", "NSE:RPOWER", "", "", "")
    50 = @("", "NSE:SAGCEM", "", "", "")
}

# New rows 33:50 did not exist before - give column A the same style as the
# existing index column (A2) before writing into them.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A33:A50").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Column A holds the 0-based row index for every data row (2-50 => 0-48).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2

    $vals = $data[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
